$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16, shifting existing rows 16:111 down to 17:112
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16 with the new data point
$ws.Range("A16").Value = 1
$ws.Range("B16").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C16").Value = "Arica y Parinacota"
$ws.Range("D16").Value = "12/06/2022"
$ws.Range("E16").Value = 15
$ws.Range("F16").Value = "Fruta"
$ws.Range("G16").Value = 100102
$ws.Range("H16").Value = "Cítricos"
$ws.Range("I16").Value = 100102005
$ws.Range("J16").Value = "Naranja"
$ws.Range("K16").Value = "Lane Late"
$ws.Range("L16").Value = "Tercera"
$ws.Range("M16").Value = 200
$ws.Range("N16").Value = 900
$ws.Range("O16").Value = 1000
$ws.Range("P16").Value = 950
$ws.Range("Q16").Value = "$/kilo (en caja de 20 kilos)"
$ws.Range("R16").Value = "Región de Coquimbo"
$ws.Range("S16").Value = 950
$ws.Range("T16").Value = 1
